# This script inserts 6 new weekly price records into the "Acelga" price
# history table (Feria Lagunitas de Puerto Montt, Los Lagos). Each new
# record is inserted at the correct position (pushing subsequent rows
# down, preserving their relative order), and the last one is appended
# right after the current last row of the table.
#
# Because the rows are processed in increasing order of the FINAL row
# index they occupy, every earlier insertion already accounts for the
# shift it causes, so by the time we reach a later target row index it is
# already valid in the (partially shifted) sheet - no extra offset
# bookkeeping is required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docena = "`$/docena de atados (4 kilos)"
$maule = "Región del Maule"

$newRows = @(
    @{RowIndex=15; Fecha=44425; Volumen=200; Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$true},
    @{RowIndex=25; Fecha=44424; Volumen=50;  Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$true},
    @{RowIndex=46; Fecha=44427; Volumen=100; Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$true},
    @{RowIndex=48; Fecha=44421; Volumen=200; Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$true},
    @{RowIndex=75; Fecha=44417; Volumen=100; Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$true},
    @{RowIndex=99; Fecha=44418; Volumen=200; Precio=4000; Unidad=$docena; Origen=$maule; PrecioKg=1000; KgUnid=4; Insert=$false}
)

foreach ($rd in $newRows) {
    $ri = $rd.RowIndex

    if ($rd.Insert) {
        $ws.Rows.Item($ri).Insert()
    }

    $ws.Range("A$ri").Value = 4
    $ws.Range("B$ri").Value = "Feria Lagunitas de Puerto Montt"
    $ws.Range("C$ri").Value = "Los Lagos"
    $ws.Range("D$ri").Value = $rd.Fecha
    $ws.Range("E$ri").Value = 10
    $ws.Range("F$ri").Value = 100112009
    $ws.Range("G$ri").Value = "Acelga"
    $ws.Range("H$ri").Value = "Sin especificar"
    $ws.Range("I$ri").Value = "Primera"
    $ws.Range("J$ri").Value = $rd.Volumen
    $ws.Range("K$ri").Value = $rd.Precio
    $ws.Range("L$ri").Value = $rd.Precio
    $ws.Range("M$ri").Value = $rd.Precio
    $ws.Range("N$ri").Value = $rd.Unidad
    $ws.Range("O$ri").Value = $rd.Origen
    $ws.Range("P$ri").Value = $rd.PrecioKg
    $ws.Range("Q$ri").Value = $rd.KgUnid
    $ws.Range("R$ri").Value = "Hortaliza"

    # Match the date-column number format used by the rest of the column.
    $ws.Range("D$ri").NumberFormat = $ws.Range("D14").NumberFormat
}
